# "Update to Entities 0.50" — retune the logistic response-curve inputs.
# Sheet1 layout: row2=M, row3=K, row4=B, row5=C (curve params) for the
# response types in columns C (Linear/Quad) and D/E/F (Logistic variants).
#
# The edit bumps the "Linear /Quad" curve's M/B/C parameters and removes
# the three "Logisttic" input columns (D:F) entirely, which drives their
# dependent formulas in rows 12:29 to recompute as 0 (or #NUM! on row 12,
# where the 0^0 power blows up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1000                   # M
$ws.Range("C4").Value = 1                      # B
$ws.Range("C5").Value = 0.55000000000000004    # C
$ws.Range("D2:F5").ClearContents()             # drop the Logisttic inputs

# Selection moves from the old C6 anchor to C4.
$ws.Range("C4").Select()
